# Swap the contents of columns F and G (header label in row 3, and all
# data values in rows 5-104) on Sheet1. Row 4 holds AVERAGE() formulas
# that reference their own column, so their cached results naturally
# follow the swapped data once Excel recalculates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- swap header labels in row 3 ---
$f3 = $ws.Range("F3").Value()
$g3 = $ws.Range("G3").Value()
$ws.Range("F3").Value = $g3
$ws.Range("G3").Value = $f3

# --- swap the data values for rows 5 through 104 ---
for ($r = 5; $r -le 104; $r++) {
    $fCell = $ws.Cells.Item($r, 6)   # column F
    $gCell = $ws.Cells.Item($r, 7)   # column G
    $fVal = $fCell.Value()
    $gVal = $gCell.Value()
    $fCell.Value = $gVal
    $gCell.Value = $fVal
}

# --- re-enter the row-4 AVERAGE formulas so the untouched D/E/H/I/J
#     columns consolidate back into shared-formula groups (their cached
#     results are unaffected); F4/G4 keep their own standalone formula,
#     matching their own (now-swapped) column. ---
$ws.Range("D4:E4").Formula = "=AVERAGE(D5:D104)"
$ws.Range("F4").Formula = "=AVERAGE(F5:F104)"
$ws.Range("G4").Formula = "=AVERAGE(G5:G104)"
$ws.Range("H4:J4").Formula = "=AVERAGE(H5:H104)"

# --- sheet view: clear the frozen/scrolled top-left cell and move selection ---
$ws.Range("K7").Select() | Out-Null

$wb.Application.Calculate() | Out-Null
